# Revert "Merge currently outstanding changes into master for important analysis"
# Removes the two stray "Straight Arrow Connector 98/99" connector shapes
# (shape ids 99 and 100) that were added to slide 1 at the top of the shape
# tree, restoring the deck to its prior state.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$s.Shapes.Item("Straight Arrow Connector 98").Delete()
$s.Shapes.Item("Straight Arrow Connector 99").Delete()
